$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.553.50'
$ws.Range("E2").Value = '  +3.52%  '
$ws.Range("D3").Value = '1.916.89'
$ws.Range("E3").Value = '  +1.89%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'248.44"
$ws.Range("E5").Value = '  +1.48%  '
$ws.Range("D6").Value = "'0.694"
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = "'44.03"
$ws.Range("E8").Value = '  +1.42%  '
$ws.Range("D9").Value = "'58.66"
$ws.Range("E9").Value = '  +9.71%  '
$ws.Range("E10").Value = '  +3.39%  '
$ws.Range("D11").Value = "'0.0764"
$ws.Range("E11").Value = '  +3.28%  '
$ws.Range("D12").Value = "'0.0997"
$ws.Range("E12").Value = '  +2.70%  '
$ws.Range("D13").Value = "'14.55"
$ws.Range("E13").Value = '  +7.95%  '
$ws.Range("E14").Value = '  +4.15%  '
$ws.Range("D15").Value = '2.200.82'
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("D16").Value = "'5.12"
$ws.Range("E16").Value = '  +4.57%  '
$ws.Range("D17").Value = '1.923.28'
$ws.Range("E17").Value = '  +1.96%  '
$ws.Range("D18").Value = '36.510.50'
$ws.Range("E18").Value = '  +3.18%  '
$ws.Range("D19").Value = "'74.13"
$ws.Range("E19").Value = '  +1.76%  '
$ws.Range("D20").Value = '0.0₃0857'
$ws.Range("E20").Value = '  +4.61%  '
$ws.Range("D21").Value = "'250.59"
$ws.Range("E21").Value = '  +2.90%  '
$ws.Range("D22").Value = "'13.21"
$ws.Range("E22").Value = '  +3.50%  '
$ws.Range("E23").Value = '  +5.18%  '
$ws.Range("D24").Value = "'2.69"
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").Value = "'2.18"
$ws.Range("E26").Value = '  +1.21%  '
$ws.Range("D27").Value = "'167.48"
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").Value = "'8.79"
$ws.Range("E28").Value = '  +2.98%  '
$ws.Range("D29").Value = "'18.72"
$ws.Range("E29").Value = '  +2.36%  '
$ws.Range("E30").Value = '  +1.49%  '
$ws.Range("E31").Value = '  +6.79%  '
$ws.Range("D32").Value = "'0.0607"
$ws.Range("E32").Value = '  +3.62%  '
$ws.Range("D33").Value = "'2.00"
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("D34").Value = "'4.33"
$ws.Range("E34").Value = '  +4.86%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  +18.66%  '
$ws.Range("D37").Value = "'1.49"
$ws.Range("E37").Value = '  -13.79%  '
$ws.Range("D38").Value = "'0.874"
$ws.Range("E38").Value = '  +3.85%  '
$ws.Range("D39").Value = "'17.63"
$ws.Range("E39").Value = '  +45.58%  '
$ws.Range("E40").Value = '  +3.63%  '
$ws.Range("D41").Value = "'106.80"
$ws.Range("E41").Value = '  +11.27%  '
$ws.Range("E42").Value = '  +5.44%  '
$ws.Range("D43").Value = "'17.12"
$ws.Range("E43").Value = '  -1.91%  '
$ws.Range("E44").Value = '  +2.92%  '
$ws.Range("D45").Value = '1.337.57'
$ws.Range("E45").Value = '  +2.59%  '
$ws.Range("D46").Value = "'2.35"
$ws.Range("E46").Value = '  +0.46%  '
$ws.Range("E47").Value = '  +4.69%  '
$ws.Range("D48").Value = "'0.0813"
$ws.Range("E48").Value = '  +1.93%  '
$ws.Range("D49").Value = "'2.79"
$ws.Range("E49").Value = '  +2.58%  '
$ws.Range("D50").Value = "'6.41"
$ws.Range("E50").Value = '  +3.18%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").Value = "'42.91"
$ws.Range("E51").Value = '  +2.11%  '
